$wb = $excel.ActiveWorkbook

# --- Rename "Attempts" sheet to "Send Attempt" ---
$ws2 = $wb.Worksheets.Item("Attempts")
$ws2.Name = "Send Attempt"

# --- Update attempt counts on the (renamed) Send Attempt sheet ---
$ws2.Cells.Item(2, 4).Value = 1
$ws2.Cells.Item(2, 6).Value = 0
$ws2.Cells.Item(2, 7).Value = 0
$ws2.Cells.Item(2, 9).Value = 0
$ws2.Cells.Item(2, 27).Value = 0
$ws2.Cells.Item(2, 41).Value = 0
$ws2.Cells.Item(3, 2).Value = 0
$ws2.Cells.Item(3, 3).Value = 0
$ws2.Cells.Item(3, 4).Value = 0
$ws2.Cells.Item(3, 5).Value = 0
$ws2.Cells.Item(3, 6).Value = 0
$ws2.Cells.Item(3, 7).Value = 0
$ws2.Cells.Item(3, 9).Value = 0
$ws2.Cells.Item(3, 41).Value = 0
$ws2.Cells.Item(4, 2).Value = 0
$ws2.Cells.Item(4, 3).Value = 0
$ws2.Cells.Item(4, 4).Value = 0
$ws2.Cells.Item(4, 8).Value = 0
$ws2.Cells.Item(4, 12).Value = 0
$ws2.Cells.Item(4, 13).Value = 0
$ws2.Cells.Item(4, 41).Value = 0
$ws2.Cells.Item(5, 3).Value = 0
$ws2.Cells.Item(5, 4).Value = 0
$ws2.Cells.Item(5, 5).Value = 0
$ws2.Cells.Item(6, 2).Value = 0
$ws2.Cells.Item(6, 3).Value = 0
$ws2.Cells.Item(6, 4).Value = 0
$ws2.Cells.Item(6, 8).Value = 0
$ws2.Cells.Item(6, 14).Value = 0
$ws2.Cells.Item(6, 41).Value = 0
$ws2.Cells.Item(7, 6).Value = 0
$ws2.Cells.Item(7, 8).Value = 0
$ws2.Cells.Item(7, 9).Value = 0
$ws2.Cells.Item(7, 13).Value = 0
$ws2.Cells.Item(7, 41).Value = 0
$ws2.Cells.Item(8, 6).Value = 0
$ws2.Cells.Item(8, 11).Value = 0
$ws2.Cells.Item(8, 12).Value = 0
$ws2.Cells.Item(8, 13).Value = 0
$ws2.Cells.Item(8, 14).Value = 0
$ws2.Cells.Item(9, 10).Value = 0
$ws2.Cells.Item(9, 12).Value = 0
$ws2.Cells.Item(9, 15).Value = 0
$ws2.Cells.Item(9, 19).Value = 0
$ws2.Cells.Item(9, 20).Value = 0
$ws2.Cells.Item(9, 41).Value = 0
$ws2.Cells.Item(10, 2).Value = 0
$ws2.Cells.Item(10, 8).Value = 0
$ws2.Cells.Item(10, 9).Value = 0
$ws2.Cells.Item(10, 11).Value = 0
$ws2.Cells.Item(10, 12).Value = 0
$ws2.Cells.Item(11, 24).Value = 0
$ws2.Cells.Item(11, 26).Value = 0
$ws2.Cells.Item(11, 28).Value = 0
$ws2.Cells.Item(11, 30).Value = 0
$ws2.Cells.Item(12, 7).Value = 0
$ws2.Cells.Item(12, 9).Value = 0
$ws2.Cells.Item(12, 35).Value = 0
$ws2.Cells.Item(12, 37).Value = 0
$ws2.Cells.Item(12, 39).Value = 0
$ws2.Cells.Item(12, 41).Value = 0
$ws2.Cells.Item(13, 17).Value = 0
$ws2.Cells.Item(13, 35).Value = 0
$ws2.Cells.Item(13, 36).Value = 0
$ws2.Cells.Item(13, 37).Value = 0
$ws2.Cells.Item(13, 38).Value = 0
$ws2.Cells.Item(14, 19).Value = 0
$ws2.Cells.Item(14, 28).Value = 0
$ws2.Cells.Item(14, 30).Value = 0
$ws2.Cells.Item(14, 31).Value = 0
$ws2.Cells.Item(14, 32).Value = 0
$ws2.Cells.Item(14, 34).Value = 0
$ws2.Cells.Item(14, 39).Value = 0
$ws2.Cells.Item(14, 41).Value = 0
$ws2.Cells.Item(15, 2).Value = 0
$ws2.Cells.Item(15, 4).Value = 0
$ws2.Cells.Item(15, 8).Value = 0
$ws2.Cells.Item(15, 10).Value = 0
$ws2.Cells.Item(15, 13).Value = 0
$ws2.Cells.Item(15, 15).Value = 0
$ws2.Cells.Item(15, 17).Value = 0
$ws2.Cells.Item(15, 33).Value = 0
$ws2.Cells.Item(15, 35).Value = 0
$ws2.Cells.Item(15, 40).Value = 0
$ws2.Cells.Item(15, 41).Value = 0
$ws2.Cells.Item(16, 8).Value = 0
$ws2.Cells.Item(16, 9).Value = 0
$ws2.Cells.Item(16, 11).Value = 0
$ws2.Cells.Item(16, 13).Value = 0
$ws2.Cells.Item(16, 14).Value = 0
$ws2.Cells.Item(17, 2).Value = 0
$ws2.Cells.Item(18, 2).Value = 0
$ws2.Cells.Item(18, 7).Value = 0
$ws2.Cells.Item(18, 9).Value = 0
$ws2.Cells.Item(19, 41).Value = 0
$ws2.Cells.Item(28, 41).Value = 0
$ws2.Cells.Item(30, 9).Value = 0
$ws2.Cells.Item(30, 41).Value = 0
$ws2.Cells.Item(41, 41).Value = 0
$ws2.Cells.Item(52, 5).Value = 0
$ws2.Cells.Item(52, 13).Value = 0

# --- Update computed values on the Scores sheet ---
$ws3 = $wb.Worksheets.Item("Scores")
$ws3.Cells.Item(2, 2).Value = 560
$ws3.Cells.Item(2, 3).Value = 80
$ws3.Cells.Item(2, 4).Value = 180
$ws3.Cells.Item(2, 5).Value = 300
$ws3.Cells.Item(3, 2).Value = 0
$ws3.Cells.Item(4, 2).Value = 490
$ws3.Cells.Item(4, 7).Value = 490
$ws3.Cells.Item(5, 2).Value = 0
$ws3.Cells.Item(6, 2).Value = 0
$ws3.Cells.Item(7, 2).Value = 0
$ws3.Cells.Item(8, 2).Value = 0
$ws3.Cells.Item(9, 2).Value = 0
$ws3.Cells.Item(10, 2).Value = 0
$ws3.Cells.Item(11, 2).Value = 0
$ws3.Cells.Item(12, 2).Value = 0
$ws3.Cells.Item(13, 2).Value = 0
$ws3.Cells.Item(14, 2).Value = 0
$ws3.Cells.Item(15, 2).Value = 0
$ws3.Cells.Item(16, 2).Value = 0
$ws3.Cells.Item(17, 2).Value = 0
$ws3.Cells.Item(18, 2).Value = 0
$ws3.Cells.Item(19, 2).Value = 0
$ws3.Cells.Item(28, 2).Value = 0
$ws3.Cells.Item(30, 2).Value = 0
$ws3.Cells.Item(41, 2).Value = 0
$ws3.Cells.Item(52, 2).Value = 0

# --- Clear stale per-route score cells on the Scores sheet ---
$ws3.Cells.Item(2, 7).ClearContents()
$ws3.Cells.Item(2, 8).ClearContents()
$ws3.Cells.Item(2, 10).ClearContents()
$ws3.Cells.Item(2, 28).ClearContents()
$ws3.Cells.Item(2, 42).ClearContents()
$ws3.Cells.Item(3, 3).ClearContents()
$ws3.Cells.Item(3, 4).ClearContents()
$ws3.Cells.Item(3, 5).ClearContents()
$ws3.Cells.Item(3, 6).ClearContents()
$ws3.Cells.Item(3, 7).ClearContents()
$ws3.Cells.Item(3, 8).ClearContents()
$ws3.Cells.Item(3, 10).ClearContents()
$ws3.Cells.Item(3, 42).ClearContents()
$ws3.Cells.Item(4, 3).ClearContents()
$ws3.Cells.Item(4, 4).ClearContents()
$ws3.Cells.Item(4, 5).ClearContents()
$ws3.Cells.Item(4, 9).ClearContents()
$ws3.Cells.Item(4, 13).ClearContents()
$ws3.Cells.Item(4, 14).ClearContents()
$ws3.Cells.Item(4, 42).ClearContents()
$ws3.Cells.Item(5, 4).ClearContents()
$ws3.Cells.Item(5, 5).ClearContents()
$ws3.Cells.Item(5, 6).ClearContents()
$ws3.Cells.Item(6, 3).ClearContents()
$ws3.Cells.Item(6, 4).ClearContents()
$ws3.Cells.Item(6, 5).ClearContents()
$ws3.Cells.Item(6, 9).ClearContents()
$ws3.Cells.Item(6, 15).ClearContents()
$ws3.Cells.Item(6, 42).ClearContents()
$ws3.Cells.Item(7, 7).ClearContents()
$ws3.Cells.Item(7, 9).ClearContents()
$ws3.Cells.Item(7, 10).ClearContents()
$ws3.Cells.Item(7, 14).ClearContents()
$ws3.Cells.Item(7, 42).ClearContents()
$ws3.Cells.Item(8, 7).ClearContents()
$ws3.Cells.Item(8, 12).ClearContents()
$ws3.Cells.Item(8, 13).ClearContents()
$ws3.Cells.Item(8, 14).ClearContents()
$ws3.Cells.Item(8, 15).ClearContents()
$ws3.Cells.Item(9, 11).ClearContents()
$ws3.Cells.Item(9, 13).ClearContents()
$ws3.Cells.Item(9, 16).ClearContents()
$ws3.Cells.Item(9, 20).ClearContents()
$ws3.Cells.Item(9, 21).ClearContents()
$ws3.Cells.Item(9, 42).ClearContents()
$ws3.Cells.Item(10, 3).ClearContents()
$ws3.Cells.Item(10, 9).ClearContents()
$ws3.Cells.Item(10, 10).ClearContents()
$ws3.Cells.Item(10, 12).ClearContents()
$ws3.Cells.Item(10, 13).ClearContents()
$ws3.Cells.Item(11, 25).ClearContents()
$ws3.Cells.Item(11, 27).ClearContents()
$ws3.Cells.Item(11, 29).ClearContents()
$ws3.Cells.Item(11, 31).ClearContents()
$ws3.Cells.Item(12, 8).ClearContents()
$ws3.Cells.Item(12, 10).ClearContents()
$ws3.Cells.Item(12, 36).ClearContents()
$ws3.Cells.Item(12, 38).ClearContents()
$ws3.Cells.Item(12, 40).ClearContents()
$ws3.Cells.Item(12, 42).ClearContents()
$ws3.Cells.Item(13, 18).ClearContents()
$ws3.Cells.Item(13, 36).ClearContents()
$ws3.Cells.Item(13, 37).ClearContents()
$ws3.Cells.Item(13, 38).ClearContents()
$ws3.Cells.Item(13, 39).ClearContents()
$ws3.Cells.Item(14, 20).ClearContents()
$ws3.Cells.Item(14, 29).ClearContents()
$ws3.Cells.Item(14, 31).ClearContents()
$ws3.Cells.Item(14, 32).ClearContents()
$ws3.Cells.Item(14, 33).ClearContents()
$ws3.Cells.Item(14, 35).ClearContents()
$ws3.Cells.Item(14, 40).ClearContents()
$ws3.Cells.Item(14, 42).ClearContents()
$ws3.Cells.Item(15, 3).ClearContents()
$ws3.Cells.Item(15, 5).ClearContents()
$ws3.Cells.Item(15, 9).ClearContents()
$ws3.Cells.Item(15, 11).ClearContents()
$ws3.Cells.Item(15, 14).ClearContents()
$ws3.Cells.Item(15, 16).ClearContents()
$ws3.Cells.Item(15, 18).ClearContents()
$ws3.Cells.Item(15, 34).ClearContents()
$ws3.Cells.Item(15, 36).ClearContents()
$ws3.Cells.Item(15, 41).ClearContents()
$ws3.Cells.Item(15, 42).ClearContents()
$ws3.Cells.Item(16, 9).ClearContents()
$ws3.Cells.Item(16, 10).ClearContents()
$ws3.Cells.Item(16, 12).ClearContents()
$ws3.Cells.Item(16, 14).ClearContents()
$ws3.Cells.Item(16, 15).ClearContents()
$ws3.Cells.Item(17, 3).ClearContents()
$ws3.Cells.Item(18, 3).ClearContents()
$ws3.Cells.Item(18, 8).ClearContents()
$ws3.Cells.Item(18, 10).ClearContents()
$ws3.Cells.Item(19, 42).ClearContents()
$ws3.Cells.Item(28, 42).ClearContents()
$ws3.Cells.Item(30, 10).ClearContents()
$ws3.Cells.Item(30, 42).ClearContents()
$ws3.Cells.Item(41, 42).ClearContents()
$ws3.Cells.Item(52, 6).ClearContents()
$ws3.Cells.Item(52, 14).ClearContents()
